$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ntf3"
$ws.Range("C2").Value = "Ntrk2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.305917
$ws.Range("H2").Value = 15.917751
$ws.Range("I2").Value = 0.4336744870332215
$ws.Range("J2").Value = 0.4336744870332215
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8220243333333334
$ws.Range("N2").Value = 2.466073
$ws.Range("O2").Value = 0.03815249372618141
$ws.Range("P2").Value = 0.03815249372618141
$ws.Range("Q2").Value = 4.361592884647
$ws.Range("R2").Value = 39.254335961823
$ws.Range("S2").Value = 0.01654576314573993
$ws.Range("T2").Value = 0.01654576314573992

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ntf3"
$ws.Range("C3").Value = "Ntrk2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.305917
$ws.Range("H3").Value = 15.917751
$ws.Range("I3").Value = 0.4336744870332215
$ws.Range("J3").Value = 0.4336744870332215
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 14.52590566666666
$ws.Range("N3").Value = 43.57771699999999
$ws.Range("O3").Value = 0.6741887099221348
$ws.Range("P3").Value = 0.6741887099221348
$ws.Range("Q3").Value = 77.07324981716299
$ws.Range("R3").Value = 693.6592483544669
$ws.Range("S3").Value = 0.2923784429390712
$ws.Range("T3").Value = 0.2923784429390712

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ntf3"
$ws.Range("C4").Value = "Ntrk2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.305917
$ws.Range("H4").Value = 15.917751
$ws.Range("I4").Value = 0.4336744870332215
$ws.Range("J4").Value = 0.4336744870332215
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.197826333333334
$ws.Range("N4").Value = 18.593479
$ws.Range("O4").Value = 0.2876587963516838
$ws.Range("P4").Value = 0.2876587963516837
$ws.Range("Q4").Value = 32.885152105081
$ws.Range("R4").Value = 295.966368945729
$ws.Range("S4").Value = 0.1247502809484104
$ws.Range("T4").Value = 0.1247502809484104

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ntf3"
$ws.Range("C5").Value = "Ntrk2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.528563666666667
$ws.Range("H5").Value = 7.585691000000001
$ws.Range("I5").Value = 0.2066699405724794
$ws.Range("J5").Value = 0.2066699405724794
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8220243333333334
$ws.Range("N5").Value = 2.466073
$ws.Range("O5").Value = 0.03815249372618141
$ws.Range("P5").Value = 0.03815249372618141
$ws.Range("Q5").Value = 2.078540862382556
$ws.Range("R5").Value = 18.706867761443
$ws.Range("S5").Value = 0.007884973611081807
$ws.Range("T5").Value = 0.007884973611081807

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ntf3"
$ws.Range("C6").Value = "Ntrk2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.528563666666667
$ws.Range("H6").Value = 7.585691000000001
$ws.Range("I6").Value = 0.2066699405724794
$ws.Range("J6").Value = 0.2066699405724794
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 14.52590566666666
$ws.Range("N6").Value = 43.57771699999999
$ws.Range("O6").Value = 0.6741887099221348
$ws.Range("P6").Value = 0.6741887099221348
$ws.Range("Q6").Value = 36.72967729416077
$ws.Range("R6").Value = 330.567095647447
$ws.Range("S6").Value = 0.1393345406142442
$ws.Range("T6").Value = 0.1393345406142442

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ntf3"
$ws.Range("C7").Value = "Ntrk2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.528563666666667
$ws.Range("H7").Value = 7.585691000000001
$ws.Range("I7").Value = 0.2066699405724794
$ws.Range("J7").Value = 0.2066699405724794
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 6.197826333333334
$ws.Range("N7").Value = 18.593479
$ws.Range("O7").Value = 0.2876587963516838
$ws.Range("P7").Value = 0.2876587963516837
$ws.Range("Q7").Value = 15.67159847877656
$ws.Range("R7").Value = 141.044386308989
$ws.Range("S7").Value = 0.05945042634715346
$ws.Range("T7").Value = 0.05945042634715345

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Ntf3"
$ws.Range("C8").Value = "Ntrk2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.400310999999999
$ws.Range("H8").Value = 13.200933
$ws.Range("I8").Value = 0.359655572394299
$ws.Range("J8").Value = 0.359655572394299
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.8220243333333334
$ws.Range("N8").Value = 2.466073
$ws.Range("O8").Value = 0.03815249372618141
$ws.Range("P8").Value = 0.03815249372618141
$ws.Range("Q8").Value = 3.617162716234333
$ws.Range("R8").Value = 32.554464446109
$ws.Range("S8").Value = 0.01372175696935968
$ws.Range("T8").Value = 0.01372175696935968

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Ntf3"
$ws.Range("C9").Value = "Ntrk2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.400310999999999
$ws.Range("H9").Value = 13.200933
$ws.Range("I9").Value = 0.359655572394299
$ws.Range("J9").Value = 0.359655572394299
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 14.52590566666666
$ws.Range("N9").Value = 43.57771699999999
$ws.Range("O9").Value = 0.6741887099221348
$ws.Range("P9").Value = 0.6741887099221348
$ws.Range("Q9").Value = 63.91850248999565
$ws.Range("R9").Value = 575.2665224099609
$ws.Range("S9").Value = 0.2424757263688194
$ws.Range("T9").Value = 0.2424757263688194

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Ntf3"
$ws.Range("C10").Value = "Ntrk2"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.400310999999999
$ws.Range("H10").Value = 13.200933
$ws.Range("I10").Value = 0.359655572394299
$ws.Range("J10").Value = 0.359655572394299
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 6.197826333333334
$ws.Range("N10").Value = 18.593479
$ws.Range("O10").Value = 0.2876587963516838
$ws.Range("P10").Value = 0.2876587963516837
$ws.Range("Q10").Value = 27.27236339065633
$ws.Range("R10").Value = 245.451270515907
$ws.Range("S10").Value = 0.1034580890561199
$ws.Range("T10").Value = 0.1034580890561199

# Remove the now-unused rows 11-13 (previously the Resolving-Mac target-cluster rows)
$ws.Rows("11:13").Delete()

Write-Output "Update complete"